$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "63.949.15"
$ws.Range("E2").Value = "  -2.82%  "

$ws.Range("D3").Value = "2.628.31"
$ws.Range("E3").Value = "  -1.07%  "

$ws.Range("E4").Value = "  +0.03%  "

$ws.Range("D5").Value = "'577.40"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -3.42%  "

$ws.Range("D6").Value = "'156.78"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.30%  "

$ws.Range("E7").Value = "  +0.10%  "

$ws.Range("D8").Value = "'0.632"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +0.42%  "

$ws.Range("D9").Value = "'0.120"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -4.64%  "

$ws.Range("D10").Value = "'5.82"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +0.32%  "

$ws.Range("D11").Value = "'0.385"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -2.84%  "

$ws.Range("E12").Value = "  -0.33%  "

$ws.Range("D13").Value = "'28.46"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -0.38%  "

$ws.Range("D14").Value = "3.108.49"
$ws.Range("E14").Value = "  -0.77%  "

$ws.Range("E15").Value = "  -6.01%  "

$ws.Range("D16").Value = "63.825.93"
$ws.Range("E16").Value = "  -2.76%  "

$ws.Range("D17").Value = "2.628.34"
$ws.Range("E17").Value = "  -0.42%  "

$ws.Range("D18").Value = "'12.15"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -3.40%  "

$ws.Range("E19").Value = "  +3.26%  "

$ws.Range("D20").Value = "'4.61"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -2.64%  "

$ws.Range("D21").Value = "'345.10"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -1.27%  "

$ws.Range("E22").Value = "  -0.17%  "

$ws.Range("E23").Value = "  -2.53%  "

$ws.Range("E24").Value = "  +1.40%  "

$ws.Range("E25").Value = "  -1.90%  "

$ws.Range("D26").Value = "'598.90"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +8.60%  "

$ws.Range("D27").Value = "'9.26"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -3.72%  "

$ws.Range("E28").Value = "  -0.28%  "

$ws.Range("E29").Value = "  -0.83%  "

$ws.Range("E30").Value = "  +0.00%  "

$ws.Range("D31").Value = "'7.91"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -0.04%  "

$ws.Range("E33").Value = "  -1.82%  "

$ws.Range("D34").Value = "'6.62"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +2.16%  "

$ws.Range("D35").Value = "'5.33"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -1.76%  "

$ws.Range("E36").Value = "  -2.10%  "

$ws.Range("D37").Value = "'19.93"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -2.07%  "

$ws.Range("D38").Value = "'0.997"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -0.18%  "

$ws.Range("D39").Value = "'154.88"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +0.05%  "

$ws.Range("E40").Value = "  -2.56%  "

$ws.Range("E41").Value = "  -0.02%  "

$ws.Range("D42").Value = "'41.55"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -2.21%  "

$ws.Range("D43").Value = "'2.43"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +6.16%  "

$ws.Range("D44").Value = "'157.26"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -2.77%  "

$ws.Range("D45").Value = "'3.95"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -2.62%  "

$ws.Range("D46").Value = "'23.22"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +2.75%  "

$ws.Range("D47").Value = "'0.0599"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -0.96%  "

$ws.Range("E48").Value = "  +2.18%  "

$ws.Range("D49").Value = "'0.630"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -1.14%  "

$ws.Range("E50").Value = "  -1.78%  "

$ws.Range("E51").Value = "  -3.25%  "
